$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9826866984367371
$ws.Range("B1").Value = 1.941970586776733
$ws.Range("C1").Value = 8.32557487487793
$ws.Range("D1").Value = 2.869608640670776
$ws.Range("E1").Value = 1.436596035957336
